$d = $word.ActiveDocument
Write-Host "Count:" $d.Styles.Count
$s4 = $d.Styles.Item(4)
Write-Host "item4:" $s4.NameLocal
$s4.Delete()
Write-Host "Count after delete4:" $d.Styles.Count

$s3 = $d.Styles.Item(3)
Write-Host "item3:" $s3.NameLocal
$s3.Delete()
Write-Host "Count after delete3:" $d.Styles.Count

$s2 = $d.Styles.Item(2)
Write-Host "item2:" $s2.NameLocal
$s2.Delete()
Write-Host "Count after delete2:" $d.Styles.Count
